# Efnb2-Pecam1 LR-pair sheet: refresh the "ECs" cluster's underlying TPM-derived
# expression numbers (ligand + receptor avg/total) and every downstream
# specificity / edge-weight value that is computed from them, per the new
# "natmiOut_TPM" script run ("update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,7).Value = 42.75280866666666
$ws.Cells.Item(2,8).Value = 128.258426
$ws.Cells.Item(2,9).Value = 0.8529286054750734
$ws.Cells.Item(2,10).Value = 0.8529286054750735
$ws.Cells.Item(2,13).Value = 206.8463543333333
$ws.Cells.Item(2,14).Value = 620.5390629999999
$ws.Cells.Item(2,15).Value = 0.9727792429746633
$ws.Cells.Item(2,16).Value = 0.9727792429746634
$ws.Cells.Item(2,17).Value = 8843.262610210535
$ws.Cells.Item(2,18).Value = 79589.36349189482
$ws.Cells.Item(2,19).Value = 0.8297112431454772
$ws.Cells.Item(2,20).Value = 0.8297112431454774

# Row 3
$ws.Cells.Item(3,7).Value = 42.75280866666666
$ws.Cells.Item(3,8).Value = 128.258426
$ws.Cells.Item(3,9).Value = 0.8529286054750734
$ws.Cells.Item(3,10).Value = 0.8529286054750735
$ws.Cells.Item(3,13).Value = 3.181559666666666
$ws.Cells.Item(3,14).Value = 9.544678999999999
$ws.Cells.Item(3,15).Value = 0.01496258038481643
$ws.Cells.Item(3,16).Value = 0.01496258038481643
$ws.Cells.Item(3,17).Value = 136.0206116905837
$ws.Cells.Item(3,18).Value = 1224.185505215254
$ws.Cells.Item(3,19).Value = 0.01276201282193016
$ws.Cells.Item(3,20).Value = 0.01276201282193017

# Row 4
$ws.Cells.Item(4,7).Value = 42.75280866666666
$ws.Cells.Item(4,8).Value = 128.258426
$ws.Cells.Item(4,9).Value = 0.8529286054750734
$ws.Cells.Item(4,10).Value = 0.8529286054750735
$ws.Cells.Item(4,15).Value = 0.01225817664052023
$ws.Cells.Item(4,16).Value = 0.01225817664052023
$ws.Cells.Item(4,17).Value = 111.4356375686895
$ws.Cells.Item(4,18).Value = 1002.920738118206
$ws.Cells.Item(4,19).Value = 0.01045534950766604
$ws.Cells.Item(4,20).Value = 0.01045534950766604

# Row 5
$ws.Cells.Item(5,9).Value = 0.04642608686423023
$ws.Cells.Item(5,10).Value = 0.04642608686423023
$ws.Cells.Item(5,13).Value = 206.8463543333333
$ws.Cells.Item(5,14).Value = 620.5390629999999
$ws.Cells.Item(5,15).Value = 0.9727792429746633
$ws.Cells.Item(5,16).Value = 0.9727792429746634
$ws.Cells.Item(5,17).Value = 481.3510479885436
$ws.Cells.Item(5,18).Value = 4332.159431896892
$ws.Cells.Item(5,19).Value = 0.04516233363406184
$ws.Cells.Item(5,20).Value = 0.04516233363406185

# Row 6
$ws.Cells.Item(6,9).Value = 0.04642608686423023
$ws.Cells.Item(6,10).Value = 0.04642608686423023
$ws.Cells.Item(6,13).Value = 3.181559666666666
$ws.Cells.Item(6,14).Value = 9.544678999999999
$ws.Cells.Item(6,15).Value = 0.01496258038481643
$ws.Cells.Item(6,16).Value = 0.01496258038481643
$ws.Cells.Item(6,17).Value = 7.403790531981778
$ws.Cells.Item(6,18).Value = 66.634114787836
$ws.Cells.Item(6,19).Value = 0.0006946540566585149
$ws.Cells.Item(6,20).Value = 0.000694654056658515

# Row 7
$ws.Cells.Item(7,9).Value = 0.04642608686423023
$ws.Cells.Item(7,10).Value = 0.04642608686423023
$ws.Cells.Item(7,15).Value = 0.01225817664052023
$ws.Cells.Item(7,16).Value = 0.01225817664052023
$ws.Cells.Item(7,19).Value = 0.0005690991735098701
$ws.Cells.Item(7,20).Value = 0.0005690991735098702

# Row 8
$ws.Cells.Item(8,7).Value = 5.044817999999999
$ws.Cells.Item(8,9).Value = 0.1006453076606963
$ws.Cells.Item(8,10).Value = 0.1006453076606963
$ws.Cells.Item(8,13).Value = 206.8463543333333
$ws.Cells.Item(8,14).Value = 620.5390629999999
$ws.Cells.Item(8,15).Value = 0.9727792429746633
$ws.Cells.Item(8,16).Value = 0.9727792429746634
$ws.Cells.Item(8,17).Value = 1043.502211575178
$ws.Cells.Item(8,18).Value = 9391.519904176599
$ws.Cells.Item(8,19).Value = 0.09790566619512422
$ws.Cells.Item(8,20).Value = 0.09790566619512424

# Row 9
$ws.Cells.Item(9,7).Value = 5.044817999999999
$ws.Cells.Item(9,9).Value = 0.1006453076606963
$ws.Cells.Item(9,10).Value = 0.1006453076606963
$ws.Cells.Item(9,13).Value = 3.181559666666666
$ws.Cells.Item(9,14).Value = 9.544678999999999
$ws.Cells.Item(9,15).Value = 0.01496258038481643
$ws.Cells.Item(9,16).Value = 0.01496258038481643
$ws.Cells.Item(9,19).Value = 0.001505913506227749
$ws.Cells.Item(9,20).Value = 0.001505913506227749

# Row 10
$ws.Cells.Item(10,7).Value = 5.044817999999999
$ws.Cells.Item(10,9).Value = 0.1006453076606963
$ws.Cells.Item(10,10).Value = 0.1006453076606963
$ws.Cells.Item(10,15).Value = 0.01225817664052023
$ws.Cells.Item(10,16).Value = 0.01225817664052023
$ws.Cells.Item(10,19).Value = 0.001233727959344319
$ws.Cells.Item(10,20).Value = 0.001233727959344319
